$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-CellText($row, $col, $old, $new) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # wdFindStop (0) for Wrap and wdReplaceOne (1) for Replace so the
    # substitution stays confined to this cell's Range and does not
    # leak into other cells that might contain the same text.
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1)
}

Replace-CellText 1 1 "448÷7=" "828÷3="
Replace-CellText 1 2 "143÷5=" "401÷8="
Replace-CellText 1 3 "279÷8=" "191÷7="
Replace-CellText 1 4 "252÷6=" "296÷8="
Replace-CellText 1 5 "462÷5=" "277÷4="

Replace-CellText 5 1 "432÷7=" "223÷5="
Replace-CellText 5 2 "405÷6=" "123÷3="
Replace-CellText 5 3 "644÷9=" "564÷9="
Replace-CellText 5 4 "225÷7=" "320÷8="
Replace-CellText 5 5 "680÷6=" "662÷9="

Replace-CellText 9 1 "750÷8=" "586÷4="
Replace-CellText 9 2 "609÷9=" "506÷3="
Replace-CellText 9 3 "470÷5=" "833÷6="
Replace-CellText 9 4 "980÷2=" "887÷4="
Replace-CellText 9 5 "311÷6=" "740÷6="

Replace-CellText 13 1 "961÷9=" "822÷4="
Replace-CellText 13 2 "225÷7=" "917÷4="
Replace-CellText 13 3 "259÷9=" "324÷2="
Replace-CellText 13 4 "804÷3=" "231÷4="
Replace-CellText 13 5 "695÷4=" "995÷5="

Replace-CellText 17 1 "883÷2=" "415÷6="
Replace-CellText 17 2 "679÷7=" "515÷3="
Replace-CellText 17 3 "120÷6=" "470÷9="
Replace-CellText 17 4 "393÷7=" "844÷6="
Replace-CellText 17 5 "872÷7=" "354÷6="

Write-Host "All replacements applied."
